# edit.ps1 -- "About Me" rewrite
#
# The target edit (per the supplied diff) does three things to the document body:
#   1. Drops the old bold/italic/underlined "About Me:" heading paragraph entirely.
#   2. Drops the trailing "---...End...---" paragraph entirely, but keeps its
#      _GoBack bookmark -- which is re-homed onto the (new) first paragraph.
#   3. The remaining four body paragraphs keep their exact text/formatting, they
#      just shift up to fill slots 1-4 (and the last of them inherits the simple
#      rFonts/lang rPr that used to live on the deleted trailing paragraph).
#
# Net effect: paragraph 1 ("About Me:") is replaced by what used to be paragraph 2
# (the "I am Manish..." intro) plus the relocated bookmark; paragraphs 2-4 shift up
# from the old paragraphs 3-5; and the old paragraph 6 disappears.
#
# Rather than try to splice/move runs piecemeal, we rewrite the whole body in one
# shot with Range.InsertXML (plain runs only -- InsertXML does not round-trip
# w:rStyle), then make a second pass re-applying the "Strong" character style to
# the designated phrases via Find + Range.Style.

$d = $word.ActiveDocument

$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">I am Manish, an analytical and solution-driven professional transitioning into Data Science with a strong foundation in real estate operations and business analytics. Currently, I am a Data Science Intern at AI-Varient, where I apply data-driven methodologies to support research, build predictive models, and contribute to impactful analytics solutions.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/></w:pPr><w:r><w:t xml:space="preserve">My experience spans real estate sales operations, business loan processing, and strategic coordination—strengthened by advanced skills in data interpretation, visualization, and performance reporting. I work with Python, SQL, Excel, and Power BI to create insightful dashboards, automate processes, and uncover meaningful patterns that drive informed decision-making.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/></w:pPr><w:r><w:t xml:space="preserve">Alongside traditional analytics, I leverage Generative AI tools like ChatGPT, Perplexity, and other AI platforms to accelerate research, optimize workflows, and enhance analytical output.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">I am actively seeking opportunities in Data Science, Data Analytics, and Business Analytics, where I can combine my technical expertise and business understanding to help organizations extract actionable insights, solve complex problems, and enable data-backed strategic decisions.</w:t></w:r></w:p>'
$d.Content.InsertXML($bodyXml) | Out-Null

if ($d.Paragraphs.Count -ne 4) {
    throw ("Expected 4 paragraphs after rewrite, got " + $d.Paragraphs.Count)
}

function Apply-Strong([int]$paraIndex, [string[]]$phrases) {
    $para = $d.Paragraphs($paraIndex)
    $paraEnd = $para.Range.End
    $cursor = $para.Range.Start
    foreach ($phrase in $phrases) {
        $rng = $d.Range($cursor, $paraEnd)
        $found = $rng.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            throw ("Could not find phrase '" + $phrase + "' in paragraph " + $paraIndex)
        }
        $rng.Style = "Strong"
        $cursor = $rng.End
    }
}

Apply-Strong 1 @('Manish', 'Data Science', 'real estate operations', 'business analytics', 'Data Science Intern at AI-Varient')
Apply-Strong 2 @('data interpretation, visualization, and performance reporting', 'Python, SQL, Excel, and Power BI')
Apply-Strong 3 @('Generative AI tools')
Apply-Strong 4 @('Data Science, Data Analytics, and Business Analytics')

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
